# DOMA-2542 Localization for Excel template (ticket_report_status_executor)
# - Remove stray spaces in "{d.tickets[i + 1]....}" -> "{d.tickets[i+1]....}" placeholders (row 3)
# - Theme cleanup: drop outerShdw effects, remap Cambria -> scheme minor font / Helvetica Neue

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "{d.tickets[i+1].categoryClassifier}"
$ws.Range("B3").Value = "{d.tickets[i+1].address}"
$ws.Range("C3").Value = "{d.tickets[i+1].processing}"
$ws.Range("D3").Value = "{d.tickets[i+1].completed}"
$ws.Range("E3").Value = "{d.tickets[i+1].canceled}"
$ws.Range("F3").Value = "{d.tickets[i+1].deferred}"
$ws.Range("G3").Value = "{d.tickets[i+1].closed}"
$ws.Range("H3").Value = "{d.tickets[i+1].new_or_reopened}"
